$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Process" column values (H2, H3) from "Search" to "Search & Typing"
$ws.Range("H2").Value = "Search & Typing"
$ws.Range("H3").Value = "Search & Typing"

# Update the active cell selection to reflect the new cursor position
$ws.Range("J14").Select()
